$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

function Set-CellText($row, $col, $newText) {
    $cell = $tbl.Cell($row, $col)
    $rng = $cell.Range
    # Exclude the trailing cell-mark / paragraph mark characters from the range
    $rng.End = $rng.End - 1
    $rng.Text = $newText
}

Set-CellText 1 1 "73÷4=18, 1"
Set-CellText 1 2 "82÷8=10, 2"
Set-CellText 1 3 "14÷9=1, 5"
Set-CellText 1 4 "35÷7=5, 0"
Set-CellText 1 5 "75÷9=8, 3"

Set-CellText 5 1 "38÷5=7, 3"
Set-CellText 5 2 "96÷2=48, 0"
Set-CellText 5 3 "68÷6=11, 2"
Set-CellText 5 4 "20÷5=4, 0"
Set-CellText 5 5 "98÷6=16, 2"

Set-CellText 9 1 "10÷5=2, 0"
Set-CellText 9 2 "14÷5=2, 4"
Set-CellText 9 3 "62÷8=7, 6"
Set-CellText 9 4 "93÷4=23, 1"
Set-CellText 9 5 "61÷5=12, 1"

Set-CellText 13 1 "96÷2=48, 0"
Set-CellText 13 2 "88÷8=11, 0"
Set-CellText 13 3 "96÷5=19, 1"
Set-CellText 13 4 "10÷5=2, 0"
Set-CellText 13 5 "58÷2=29, 0"

Set-CellText 17 1 "89÷2=44, 1"
Set-CellText 17 2 "21÷2=10, 1"
Set-CellText 17 3 "13÷3=4, 1"
Set-CellText 17 4 "51÷3=17, 0"
Set-CellText 17 5 "90÷6=15, 0"
